# Add a "metadata" worksheet (after the existing "data" sheet) that records
# the PanelApp query provenance for this export: which panel, which
# version, and when/how it was fetched.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# New sheet goes right after "data", matching the authored sheet order.
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$ws.Name = "metadata"

# Header row (B1:G1) - data_name / data_id / ... columns. A1 is left blank.
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Reuse the existing bold/bordered header style from the "data" sheet
# instead of minting a fresh (visually-identical) style entry.
$dataSheet.Range("B1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)   # xlPasteFormats

# Data row 2.
$ws.Range("A2").Value = 0
$dataSheet.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)      # xlPasteFormats (same header style)

# Force D2 to stay text ("1.0") instead of being coerced to the number 1,
# then drop back to the default/unstyled cell format (no explicit style).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("B2").Value = "Undiagnosed monogenic disorder seen in a specialist genetics clinic"
$ws.Range("C2").Value = 216
$ws.Range("D2").Value = "1.0"
$ws.Range("E2").Value = "2019-10-07T08:14:58.390771Z"
$ws.Range("F2").Value = "2021-10-05 14:23:02.967568"
$ws.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/216/?format=json"
$ws.Range("D2").Style = "Normal"

# Leave "data" as the active sheet/selection, as it was before this edit.
$dataSheet.Activate()
[void]$dataSheet.Range("A1").Select()
